$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51; existing rows 51-97 shift down to 52-98.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new record's data.
$ws.Range("A51").Value = 6
$ws.Range("B51").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C51").Value = "Metropolitana"
$ws.Range("D51").Value = 44790
$ws.Range("E51").Value = 13
$ws.Range("F51").Value = 100114007
$ws.Range("G51").Value = "Jengibre"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 200
$ws.Range("K51").Value = 11000
$ws.Range("L51").Value = 12000
$ws.Range("M51").Value = 11600
$ws.Range("N51").Value = "$/caja 13 kilos"
$ws.Range("O51").Value = "Perú"
$ws.Range("P51").Value = 892
$ws.Range("Q51").Value = 13
$ws.Range("R51").Value = "Hortaliza"
